# Applies updated CPU process statistics to the procesosCPU sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (nginx)
$ws.Cells.Item(2, 2).Value = 8151
$ws.Cells.Item(2, 3).Value = "'1.73%"

# Row 3 (app)
$ws.Cells.Item(3, 2).Value = 8151
$ws.Cells.Item(3, 3).Value = "'9.15%"
$ws.Cells.Item(3, 4).Value = "'0.03%"

# Row 4 (broker-broker_server-1)
$ws.Cells.Item(4, 2).Value = 8097
$ws.Cells.Item(4, 3).Value = "'9.98%"

# Row 5 (bara_text_recognition)
$ws.Cells.Item(5, 2).Value = 8055
$ws.Cells.Item(5, 3).Value = "'99.90%"

# Row 6 (bara_tag_detection)
$ws.Cells.Item(6, 2).Value = 8055
$ws.Cells.Item(6, 3).Value = "'99.76%"

# Row 7 (bara_product_list_checker)
$ws.Cells.Item(7, 2).Value = 8055
$ws.Cells.Item(7, 3).Value = "'2.50%"

# Row 8 (bara_text_detection)
$ws.Cells.Item(8, 2).Value = 8055
$ws.Cells.Item(8, 3).Value = "'99.98%"

# Row 9 (bara_minio)
$ws.Cells.Item(9, 2).Value = 8148
$ws.Cells.Item(9, 3).Value = "'9.99%"

# Row 10 (bara_rabbitmq)
$ws.Cells.Item(10, 2).Value = 8100
$ws.Cells.Item(10, 3).Value = "'94.87%"
$ws.Cells.Item(10, 4).Value = "'0.00%"

# Row 11 (bara_redis)
$ws.Cells.Item(11, 2).Value = 8148
$ws.Cells.Item(11, 3).Value = "'9.86%"
$ws.Cells.Item(11, 4).Value = "'0.12%"

# Row 12 (bara_database_2)
$ws.Cells.Item(12, 2).Value = 8150
$ws.Cells.Item(12, 3).Value = "'9.82%"
$ws.Cells.Item(12, 4).Value = "'0.27%"

# Row 13 (bara_database)
$ws.Cells.Item(13, 2).Value = 176
$ws.Cells.Item(13, 3).Value = "'3.79%"
$ws.Cells.Item(13, 4).Value = "'0.54%"
